# Add a new weekly price observation for "Arándano (blue)" at the top of the
# date-ordered block (row 241), pushing the existing rows down by one.
#
# The new row is created by copying the row that currently sits at 241 (so it
# inherits the unchanged descriptive columns: Mercado, Región, Codreg, Tipo,
# Producto, Categoría, Variedad, Calidad, Unidad de comercialización, Kg/unidad)
# and inserting that copy above it. Afterwards we overwrite only the columns
# that actually carry new data for this observation: Fecha, Volumen, Precio
# mínimo/máximo/promedio ponderado, Origen and Precio $/Kg.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(241).Copy()
$ws.Rows.Item(241).Insert()

$ws.Range("D241").Value = 45229
$ws.Range("M241").Value = 280
$ws.Range("N241").Value = 10000
$ws.Range("O241").Value = 10000
$ws.Range("P241").Value = 10000
$ws.Range("R241").Value = "Provincia de Curicó"
$ws.Range("S241").Value = 5000
